$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '35.570.98'
$ws.Range('E2').Value = '  +1.48%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.919.03'
$ws.Range('E3').Value = '  +3.56%  '
$ws.Range('E4').Value = '  +0.53%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '247.13'
$ws.Range('E5').Value = '  +4.19%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.652'
$ws.Range('E6').Value = '  +4.91%  '
$ws.Range('E7').Value = '  +0.44%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '42.34'
$ws.Range('E8').Value = '  +0.58%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.348'
$ws.Range('E9').Value = '  +6.23%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '49.19'
$ws.Range('E10').Value = '  +5.28%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0722'
$ws.Range('E11').Value = '  +3.98%  '
$ws.Range('E12').Value = '  +1.15%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '2.196.66'
$ws.Range('E13').Value = '  +3.51%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '12.33'
$ws.Range('E14').Value = '  +8.10%  '
$ws.Range('B15').Value = 'Polygon'
$ws.Range('C15').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.703'
$ws.Range('E15').Value = '  +3.75%  '
$ws.Range('B16').Value = 'WrappedEther'
$ws.Range('C16').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.917.64'
$ws.Range('E16').Value = '  +2.97%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '4.92'
$ws.Range('E17').Value = '  +4.83%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '35.637.17'
$ws.Range('E18').Value = '  +1.70%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '72.57'
$ws.Range('E19').Value = '  +3.29%  '
$ws.Range('D20').Value = '0.0₃0825'
$ws.Range('E20').Value = '  +3.76%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '246.34'
$ws.Range('E21').Value = '  +2.47%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '12.75'
$ws.Range('E22').Value = '  +4.77%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.87'
$ws.Range('E23').Value = '  +1.85%  '
$ws.Range('E24').Value = '  +0.48%  '
$ws.Range('E25').Value = '  +1.89%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.22'
$ws.Range('E26').Value = '  +16.75%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '171.54'
$ws.Range('E27').Value = '  +0.50%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '8.52'
$ws.Range('E28').Value = '  +6.95%  '
$ws.Range('E29').Value = '  +5.28%  '
$ws.Range('E30').Value = '  +3.00%  '
$ws.Range('E31').Value = '  +4.39%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.0574'
$ws.Range('E32').Value = '  +3.04%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.949'
$ws.Range('E33').Value = '  +22.23%  '
$ws.Range('E34').Value = '  +0.47%  '
$ws.Range('E35').Value = '  +4.16%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.75'
$ws.Range('E36').Value = '  +6.36%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.05'
$ws.Range('E37').Value = '  +1.67%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.35'
$ws.Range('E38').Value = '  +1.60%  '
$ws.Range('E39').Value = '  +3.42%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0211'
$ws.Range('E40').Value = '  +4.46%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '92.28'
$ws.Range('E41').Value = '  +0.93%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0638'
$ws.Range('E42').Value = '  +16.10%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '15.83'
$ws.Range('E43').Value = '  +7.35%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.358.00'
$ws.Range('E44').Value = '  +0.55%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.43'
$ws.Range('E45').Value = '  +3.11%  '
$ws.Range('B46').Value = 'MultiversX'
$ws.Range('C46').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '46.88'
$ws.Range('E46').Value = '  +36.20%  '
$ws.Range('B47').Value = 'Gas'
$ws.Range('C47').Value = 'https://coinranking.com/coin/hfw0nnnLtSFc7+gas-gas'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '12.69'
$ws.Range('E47').Value = '  +1.02%  '
$ws.Range('B48').Value = 'MXToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.80'
$ws.Range('E48').Value = '  +2.82%  '
$ws.Range('B49').Value = 'HuobiToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.40'
$ws.Range('E49').Value = '  -0.32%  '
$ws.Range('E50').Value = '  +0.96%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.107.87'
$ws.Range('E51').Value = '  +3.57%  '
